$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto data refresh.
# Price cells are forced to text via NumberFormat "@" (to preserve values like "1.744.43"
# as text rather than being auto-parsed as numbers/dates), then the style is reset back
# to "Normal" so no extra formatting is left applied to the cell.

$cellD2 = $ws.Range("D2")
$cellD2.NumberFormat = "@"
$cellD2.Value = "27.518.54"
$cellD2.Style = "Normal"

$cellD3 = $ws.Range("D3")
$cellD3.NumberFormat = "@"
$cellD3.Value = "1.744.76"
$cellD3.Style = "Normal"
$ws.Range("E3").Value = "  -0.33%  "

$cellD5 = $ws.Range("D5")
$cellD5.NumberFormat = "@"
$cellD5.Value = "322.75"
$cellD5.Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  +0.04%  "

$cellD7 = $ws.Range("D7")
$cellD7.NumberFormat = "@"
$cellD7.Value = "0.4444"
$cellD7.Style = "Normal"
$ws.Range("E7").Value = "  +4.49%  "

$cellD8 = $ws.Range("D8")
$cellD8.NumberFormat = "@"
$cellD8.Value = "0.3521"
$cellD8.Style = "Normal"
$ws.Range("E8").Value = "  -2.23%  "

$cellD9 = $ws.Range("D9")
$cellD9.NumberFormat = "@"
$cellD9.Value = "0.07413"
$cellD9.Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "

$cellD10 = $ws.Range("D10")
$cellD10.NumberFormat = "@"
$cellD10.Value = "41.59"
$cellD10.Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "

$cellD11 = $ws.Range("D11")
$cellD11.NumberFormat = "@"
$cellD11.Value = "1.080"
$cellD11.Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("E12").Value = "  +0.05%  "

$cellD14 = $ws.Range("D14")
$cellD14.NumberFormat = "@"
$cellD14.Value = "5.912"
$cellD14.Style = "Normal"
$ws.Range("E14").Value = "  -1.86%  "

$cellD15 = $ws.Range("D15")
$cellD15.NumberFormat = "@"
$cellD15.Value = "7.092"
$cellD15.Style = "Normal"
$ws.Range("E15").Value = "  -1.53%  "

$cellD16 = $ws.Range("D16")
$cellD16.NumberFormat = "@"
$cellD16.Value = "1.741.35"
$cellD16.Style = "Normal"
$ws.Range("E16").Value = "  -0.27%  "

$cellD17 = $ws.Range("D17")
$cellD17.NumberFormat = "@"
$cellD17.Value = "91.57"
$cellD17.Style = "Normal"
$ws.Range("E17").Value = "  -1.99%  "

$cellD19 = $ws.Range("D19")
$cellD19.NumberFormat = "@"
$cellD19.Value = "0.06386"
$cellD19.Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "

$cellD20 = $ws.Range("D20")
$cellD20.NumberFormat = "@"
$cellD20.Value = "1.000"
$cellD20.Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("E21").Value = "  -0.81%  "

$cellD22 = $ws.Range("D22")
$cellD22.NumberFormat = "@"
$cellD22.Value = "5.725"
$cellD22.Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "

$cellD23 = $ws.Range("D23")
$cellD23.NumberFormat = "@"
$cellD23.Value = "27.553.93"
$cellD23.Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

$cellD24 = $ws.Range("D24")
$cellD24.NumberFormat = "@"
$cellD24.Value = "11.11"
$cellD24.Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("E25").Value = "  +0.25%  "

$cellD26 = $ws.Range("D26")
$cellD26.NumberFormat = "@"
$cellD26.Value = "160.74"
$cellD26.Style = "Normal"
$ws.Range("E26").Value = "  -0.75%  "

$cellD27 = $ws.Range("D27")
$cellD27.NumberFormat = "@"
$cellD27.Value = "20.03"
$cellD27.Style = "Normal"
$ws.Range("E27").Value = "  -0.86%  "

$cellD28 = $ws.Range("D28")
$cellD28.NumberFormat = "@"
$cellD28.Value = "1.942.40"
$cellD28.Style = "Normal"
$ws.Range("E28").Value = "  -0.12%  "

$cellD29 = $ws.Range("D29")
$cellD29.NumberFormat = "@"
$cellD29.Value = "125.33"
$cellD29.Style = "Normal"
$ws.Range("E29").Value = "  +1.20%  "

$cellD30 = $ws.Range("D30")
$cellD30.NumberFormat = "@"
$cellD30.Value = "2.042"
$cellD30.Style = "Normal"
$ws.Range("E30").Value = "  -4.43%  "

$cellD31 = $ws.Range("D31")
$cellD31.NumberFormat = "@"
$cellD31.Value = "1.049"
$cellD31.Style = "Normal"
$ws.Range("E31").Value = "  -4.78%  "

$cellD32 = $ws.Range("D32")
$cellD32.NumberFormat = "@"
$cellD32.Value = "0.09085"
$cellD32.Style = "Normal"
$ws.Range("E32").Value = "  +2.23%  "

$cellD33 = $ws.Range("D33")
$cellD33.NumberFormat = "@"
$cellD33.Value = "3.651"
$cellD33.Style = "Normal"
$ws.Range("E33").Value = "  -0.32%  "

$cellD34 = $ws.Range("D34")
$cellD34.NumberFormat = "@"
$cellD34.Value = "5.382"
$cellD34.Style = "Normal"
$ws.Range("E34").Value = "  -3.19%  "

$cellD35 = $ws.Range("D35")
$cellD35.NumberFormat = "@"
$cellD35.Value = "0.02276"
$cellD35.Style = "Normal"
$ws.Range("E35").Value = "  -0.68%  "

$cellD36 = $ws.Range("D36")
$cellD36.NumberFormat = "@"
$cellD36.Value = "11.61"
$cellD36.Style = "Normal"
$ws.Range("E36").Value = "  -4.98%  "

$cellD37 = $ws.Range("D37")
$cellD37.NumberFormat = "@"
$cellD37.Value = "0.06037"
$cellD37.Style = "Normal"
$ws.Range("E37").Value = "  +0.64%  "

$cellD38 = $ws.Range("D38")
$cellD38.NumberFormat = "@"
$cellD38.Value = "0.2066"
$cellD38.Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "

$cellD39 = $ws.Range("D39")
$cellD39.NumberFormat = "@"
$cellD39.Value = "4.894"
$cellD39.Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "

$cellD40 = $ws.Range("D40")
$cellD40.NumberFormat = "@"
$cellD40.Value = "0.6238"
$cellD40.Style = "Normal"
$ws.Range("E40").Value = "  -1.24%  "

$cellD41 = $ws.Range("D41")
$cellD41.NumberFormat = "@"
$cellD41.Value = "1.185"
$cellD41.Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "

$cellD42 = $ws.Range("D42")
$cellD42.NumberFormat = "@"
$cellD42.Value = "1.376"
$cellD42.Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "

$cellD43 = $ws.Range("D43")
$cellD43.NumberFormat = "@"
$cellD43.Value = "7.717"
$cellD43.Style = "Normal"
$ws.Range("E43").Value = "  -2.60%  "

$cellD44 = $ws.Range("D44")
$cellD44.NumberFormat = "@"
$cellD44.Value = "13.11"
$cellD44.Style = "Normal"
$ws.Range("E44").Value = "  -2.47%  "

$ws.Range("E45").Value = "  +0.23%  "

$cellD46 = $ws.Range("D46")
$cellD46.NumberFormat = "@"
$cellD46.Value = "0.5800"
$cellD46.Style = "Normal"
$ws.Range("E46").Value = "  -1.23%  "

$ws.Range("E47").Value = "  -0.87%  "

$cellD48 = $ws.Range("D48")
$cellD48.NumberFormat = "@"
$cellD48.Value = "1.926"
$cellD48.Style = "Normal"
$ws.Range("E48").Value = "  -2.00%  "

$cellD49 = $ws.Range("D49")
$cellD49.NumberFormat = "@"
$cellD49.Value = "0.06846"
$cellD49.Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "

$cellD50 = $ws.Range("D50")
$cellD50.NumberFormat = "@"
$cellD50.Value = "1.117"
$cellD50.Style = "Normal"
$ws.Range("E50").Value = "  -3.85%  "

$cellD51 = $ws.Range("D51")
$cellD51.NumberFormat = "@"
$cellD51.Value = "71.45"
$cellD51.Style = "Normal"
$ws.Range("E51").Value = "  -2.19%  "
